# Insert a new data row at row 101 (pushes the existing rows 101-124 down
# to 102-125) and populate it with the new weekly price observation.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(101).Insert()

# Columns that are constant across every "Poroto granado" row in this sheet
# (market id/name, region code, category id/name, variety, quality, unit,
# origin, kg-per-unit, classification) - copy them from the row right below
# (which still holds the same constant values) so the new row matches the
# rest of the table.
$ws.Cells.Item(101, 1).Value = $ws.Cells.Item(102, 1).Value()
$ws.Cells.Item(101, 2).Value = $ws.Cells.Item(102, 2).Value()
$ws.Cells.Item(101, 3).Value = $ws.Cells.Item(102, 3).Value()
$ws.Cells.Item(101, 5).Value = $ws.Cells.Item(102, 5).Value()
$ws.Cells.Item(101, 6).Value = $ws.Cells.Item(102, 6).Value()
$ws.Cells.Item(101, 7).Value = $ws.Cells.Item(102, 7).Value()
$ws.Cells.Item(101, 8).Value = $ws.Cells.Item(102, 8).Value()
$ws.Cells.Item(101, 9).Value = $ws.Cells.Item(102, 9).Value()
$ws.Cells.Item(101, 14).Value = $ws.Cells.Item(102, 14).Value()
$ws.Cells.Item(101, 15).Value = $ws.Cells.Item(102, 15).Value()
$ws.Cells.Item(101, 17).Value = $ws.Cells.Item(102, 17).Value()
$ws.Cells.Item(101, 18).Value = $ws.Cells.Item(102, 18).Value()

# New values for this observation (week of 2023-03-13, serial 45015).
$ws.Cells.Item(101, 4).Value = 45015
$ws.Cells.Item(101, 10).Value = 80
$ws.Cells.Item(101, 11).Value = 28000
$ws.Cells.Item(101, 12).Value = 30000
$ws.Cells.Item(101, 13).Value = 29000
$ws.Cells.Item(101, 16).Value = 1160
